# actualizacion de asistencias al 14 de oct
# Adds two new attendance-taking columns (AA/AB/AC) to the "asistencia"
# sheet: a new date (26/09) in AA with an "obs" separator in AB, and
# another new date (03/10) in AC, mirroring the existing V/W/X/Y/Z layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asistencia")

# --- Header row (row 2): new date columns ---
$ws.Range("AA2").Value = 45926
$ws.Range("AA2").NumberFormat = "dd/mm/yy"
$ws.Range("AB2").Value = "obs"
$ws.Range("AC2").Value = 45933
$ws.Range("AC2").NumberFormat = "dd/mm/yy"

# --- Per-student attendance marks (rows 3-16) for the two new dates ---
# AA column (26/09) - some students marked late ("T") with an arrival
# time recorded in AB.
$aa = @{
    3  = "T"
    4  = "P"
    5  = "A"
    6  = "A"
    7  = "P"
    8  = "T"
    9  = "A"
    10 = "T"
    11 = "P"
    12 = "A"
    13 = "P"
    14 = "A"
    15 = "T"
    16 = "A"
}
$abTime = @{
    3  = 0.645833333333333
    8  = 0.645833333333333
    10 = 0.645833333333333
    15 = 0.645833333333333
}

# AC column (03/10) attendance marks.
$ac = @{
    3  = "A"
    4  = "P"
    5  = "A"
    6  = "A"
    7  = "A"
    8  = "A"
    9  = "A"
    10 = "A"
    11 = "A"
    12 = "A"
    13 = "A"
    14 = "A"
    15 = "A"
    16 = "A"
}

foreach ($row in 3..16) {
    $ws.Range("AA$row").Value = $aa[$row]
    if ($abTime.ContainsKey($row)) {
        $ws.Range("AB$row").Value = $abTime[$row]
        $ws.Range("AB$row").NumberFormat = "hh:mm:ss"
    }
    $ws.Range("AC$row").Value = $ac[$row]
}

# --- Summary formulas (rows 17-20) for the new columns ---
$ws.Range("AA17").Formula = '=COUNTIF(AA3:AA16,"P")'
$ws.Range("AA18").Formula = '=COUNTIF(AA3:AA16,"M")'
$ws.Range("AA19").Formula = '=COUNTIF(AA3:AA16,"T")'
$ws.Range("AA20").Formula = '=SUM(AA17:AA19)'

$ws.Range("AC17").Formula = '=COUNTIF(AC3:AC16,"P")'
$ws.Range("AC18").Formula = '=COUNTIF(AC3:AC16,"M")'
$ws.Range("AC19").Formula = '=COUNTIF(AC3:AC16,"T")'
$ws.Range("AC20").Formula = '=SUM(AC17:AC19)'

# --- Move the selection where the author left off ---
$ws.Range("AC8").Select() | Out-Null
